$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the 11/25/23 meeting: mark attendance "Y" for both SJ and GW
$ws.Range("B5").Value = "Y"
$ws.Range("C5").Value = "Y"

# Move the active selection to C5
$ws.Range("C5").Select()
